$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 209, shifting existing rows 209-293 down to 210-294.
$ws.Rows("209").Insert()

$ws.Cells.Item(209, 1).Value = 10
$ws.Cells.Item(209, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(209, 3).Value = "La Araucanía"
$ws.Cells.Item(209, 4).Value = 44636
$ws.Cells.Item(209, 5).Value = 9
$ws.Cells.Item(209, 6).Value = 100112009
$ws.Cells.Item(209, 7).Value = "Acelga"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 40
$ws.Cells.Item(209, 11).Value = 8000
$ws.Cells.Item(209, 12).Value = 8000
$ws.Cells.Item(209, 13).Value = 8000
$ws.Cells.Item(209, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(209, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(209, 16).Value = 667
$ws.Cells.Item(209, 17).Value = 12
$ws.Cells.Item(209, 18).Value = "Hortaliza"
